# "updated main GSC export data"
#
# The "Chart" sheet holds a rolling 89-day window of GSC export data
# (Date | Invalid | Valid), one row per day, oldest day first. This
# commit rolls the window forward by one day: the oldest day
# (2025-10-28, row 2) drops off, every remaining row shifts up by one,
# and a new day (2026-01-25) is appended at the bottom with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the oldest day (row 2) - this shifts all subsequent rows up by
# one, which is exactly the "rolling window" update the export does
# every day.
$ws.Rows.Item(2).Delete()

# After the delete, the last populated row is 89 (2026-01-24). Append
# the new day as the new row 90.
$lastRow = 90

# Force the date to be stored as text (matching every other date cell
# in the column) instead of letting Excel auto-convert the
# "yyyy-MM-dd"-looking string into a date serial number.
$dateCell = $ws.Range("A" + $lastRow)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-01-25"
$dateCell.ClearFormats()

$ws.Range("B" + $lastRow).Value = 0
$ws.Range("C" + $lastRow).Value = 24
